$d = $word.ActiveDocument

# Remove the legacy "_GoBack" bookmark from the "Tener dos roles" paragraph;
# it will be re-added, alone, in its own paragraph at the end of this block.
$bm = $d.Bookmarks.Item("_GoBack")
$anchorPara = $bm.Range.Paragraphs.Item(1).Index
$bm.Delete()

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Content (without the <w:p ...> wrapper) for each paragraph to insert, in
# document order, right after "Tener dos roles".
$newParagraphs = @(
    '',
    '<w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Planteamiento del problema</w:t></w:r>',
    '<w:r><w:t>Lleva preguntas</w:t></w:r>',
    '<w:r><w:t xml:space="preserve">Plantea la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>problemica</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> del proyecto</w:t></w:r>',
    '<w:r><w:t xml:space="preserve">Inicia muy general: ejemplo en el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cetis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> actualmente se generan reportes</w:t></w:r>',
    '<w:r><w:t xml:space="preserve">Se va planteando el problema conforme se desarrolla la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pagina</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> iniciando desde algo general hasta llegar a la escuela</w:t></w:r>',
    '<w:r><w:t>Las preguntas tienen que llevar a la problemática que es lo que da pie a nuestro sistema</w:t></w:r>',
    ''
)

$insertAfterIndex = $anchorPara
foreach ($content in $newParagraphs) {
    $afterPara = $d.Paragraphs.Item($insertAfterIndex)
    $r = $afterPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $insertAfterIndex = $insertAfterIndex + 1
    $newPara = $d.Paragraphs.Item($insertAfterIndex)
    $newPara.Range.InsertXML("<w:p $wNs>$content</w:p>")
}

# Re-add the bookmark, alone, as its own trailing paragraph.
$afterPara = $d.Paragraphs.Item($insertAfterIndex)
$r = $afterPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$bmPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bmPara.Range.InsertXML("<w:p $wNs><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>")

Write-Host "Done"
